$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new headers in row 1 for columns I and J
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Copy the style from the existing header cell (H1) so the new headers match formatting
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Fill in the numeric data for rows 2-4
$ws.Range("I2").Value = 8
$ws.Range("J2").Value = 8

$ws.Range("I3").Value = 8
$ws.Range("J3").Value = 9

$ws.Range("I4").Value = 5
$ws.Range("J4").Value = 5
